$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pages")
Write-Output "before:"
Write-Output $ws.AutoFilter.Range.Address()
$ws.Rows.Item(61).Insert()
Write-Output "after:"
Write-Output $ws.AutoFilter.Range.Address()
